$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("K6").Value = 2.25
$ws.Range("W6").Value = 10
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 7
$ws.Range("BA6").Value = 51
$ws.Range("BB6").Value = 126
$ws.Range("J7").Value = 3.1
$ws.Range("U7").Value = 1.88
$ws.Range("W7").Value = 6.5
$ws.Range("AH7").Value = 7.8
$ws.Range("AK7").Value = 40
$ws.Range("AN7").Value = 4.1
$ws.Range("AR7").Value = 120
$ws.Range("AU7").Value = 7.6
$ws.Range("AY7").Value = 28
$ws.Range("N8").Value = 7.1
$ws.Range("G9").Value = 1.4
$ws.Range("H9").Value = 4.4
$ws.Range("I9").Value = 5.8
$ws.Range("J9").Value = 1.85
$ws.Range("K9").Value = 2.45
$ws.Range("L9").Value = 5.5
$ws.Range("P9").Value = 4.7
$ws.Range("Q9").Value = 1.47
$ws.Range("R9").Value = 2.35
$ws.Range("U9").Value = 1.66
$ws.Range("V9").Value = 2.16
$ws.Range("W9").Value = 7.6
$ws.Range("X9").Value = 6.8
$ws.Range("Z9").Value = 8.5
$ws.Range("AA9").Value = 9
$ws.Range("AB9").Value = 17
$ws.Range("AC9").Value = 15.5
$ws.Range("AD9").Value = 7.9
$ws.Range("AE9").Value = 13.5
$ws.Range("AF9").Value = 45
$ws.Range("AH9").Value = 16.5
$ws.Range("AI9").Value = 32
$ws.Range("AJ9").Value = 15.5
$ws.Range("AK9").Value = 90
$ws.Range("AL9").Value = 45
$ws.Range("AM9").Value = 37
$ws.Range("AN9").Value = 3.4
$ws.Range("AO9").Value = 6.4
$ws.Range("AQ9").Value = 17
$ws.Range("AR9").Value = 37
$ws.Range("AU9").Value = 7.5
$ws.Range("AV9").Value = 60
$ws.Range("AW9").Value = 7.6
$ws.Range("AX9").Value = 32
$ws.Range("AY9").Value = 32
$ws.Range("AZ9").Value = 200
$ws.Range("BA9").Value = 200
$ws.Range("BB9").Value = 350
$ws.Range("G11").Value = 18.5
$ws.Range("I11").Value = 1.07
$ws.Range("J11").Value = 13.5
$ws.Range("K11").Value = 3.25
$ws.Range("L11").Value = 1.33
$ws.Range("O11").Value = 1.03
$ws.Range("P11").Value = 9
$ws.Range("Q11").Value = 1.18
$ws.Range("R11").Value = 3.84
$ws.Range("S11").Value = 1.12
$ws.Range("T11").Value = 5.3
$ws.Range("U11").Value = 2.51
$ws.Range("V11").Value = 1.5
$ws.Range("Y11").Value = 65
$ws.Range("Z11").Value = 101
$ws.Range("AB11").Value = 200
$ws.Range("AC11").Value = 21
$ws.Range("AE11").Value = 40
$ws.Range("AG11").Value = 101
$ws.Range("AH11").Value = 9.25
$ws.Range("AI11").Value = 5.8
$ws.Range("AJ11").Value = 11.5
$ws.Range("AK11").Value = 5.2
$ws.Range("AL11").Value = 10.25
$ws.Range("AM11").Value = 37
$ws.Range("AN11").Value = 18.5
$ws.Range("AO11").Value = 150
$ws.Range("AP11").Value = 100
$ws.Range("AT11").Value = 4.3
$ws.Range("AU11").Value = 13
$ws.Range("AV11").Value = 150
$ws.Range("AW11").Value = 3.05
$ws.Range("AX11").Value = 3.95
$ws.Range("AY11").Value = 17
$ws.Range("BA11").Value = 32
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 1.65
$ws.Range("L12").Value = 2.25
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 8
$ws.Range("O12").Value = 1.33
$ws.Range("P12").Value = 3.25
$ws.Range("Q12").Value = 2.05
$ws.Range("R12").Value = 1.75
$ws.Range("X12").Value = 26
$ws.Range("AC12").Value = 8.5
$ws.Range("AE12").Value = 21
$ws.Range("AI12").Value = 7
$ws.Range("AK12").Value = 12
$ws.Range("AN12").Value = 7
$ws.Range("AS12").Value = 450
$ws.Range("AU12").Value = 9.5
$ws.Range("AW12").Value = 3.5
$ws.Range("H13").Value = 3.25
$ws.Range("J13").Value = 3.2
$ws.Range("L13").Value = 3.5
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("S13").Value = 1.44
$ws.Range("T13").Value = 2.63
$ws.Range("U13").Value = 1.8
$ws.Range("V13").Value = 1.91
$ws.Range("W13").Value = 8
$ws.Range("AC13").Value = 9.5
$ws.Range("AE13").Value = 15
$ws.Range("AF13").Value = 51
$ws.Range("AG13").Value = 251
$ws.Range("AH13").Value = 9
$ws.Range("AI13").Value = 13
$ws.Range("AM13").Value = 34
$ws.Range("AO13").Value = 15
$ws.Range("AT13").Value = 2.63
$ws.Range("AY13").Value = 26
$ws.Range("BA13").Value = 81
$ws.Range("BB13").Value = 201
$ws.Range("Q14").Value = 1.95
$ws.Range("R14").Value = 1.9
